# bulk_user_details.xlsx — swap in a fresh set of iAuthor test-case
# credentials for row 2 and drop the extra candidate row (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second test-case row entirely; remaining rows shift up.
$ws.Range("A3:H3").Delete() | Out-Null

# New credentials/names for the remaining test-case row.
$ws.Range("A2").Value = "axkhA153"
$ws.Range("B2").Value = 231016173
$ws.Range("C2").Value = "nwphlqq84"
$ws.Range("D2").Value = "yu#&KV68"
$ws.Range("F2").Value = "wmbavLmo"
$ws.Range("G2").Value = "BZsr"

# Keep the selection in sync with the now-smaller used range.
$ws.Range("A1:H2").Select() | Out-Null
